# BankManagerSuite.xlsx: "Data Driven Extended With Docker"
#
# - TestCases!B3 and TestData!A4 held the literal "N" (Runmode=No);
#   flip them to "Y" so both suites now run. This leaves the shared
#   string "N" unreferenced, so Excel drops it from sharedStrings.xml
#   on save (re-indexing every subsequent <v> automatically).
# - Sheet3 was a blank placeholder tab; remove it.
# - Refresh the active-cell selections on TestCases/TestData and drop
#   the stale topLeftCell scroll position on TestData now that row 11
#   is no longer the focal row.

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestData  = $wb.Worksheets.Item("TestData")

# Flip the two "N" (no-run) flags to "Y" (run).
$wsTestCases.Range("B3").Value = "Y"
$wsTestData.Range("A4").Value = "Y"

# Drop the now-empty placeholder sheet.
$wsSheet3 = $wb.Worksheets.Item("Sheet3")
$wsSheet3.Delete()

# Update on-screen selections to match the post-edit focal cells.
$wsTestCases.Activate()
$wsTestCases.Range("B3").Select()

$wsTestData.Activate()
$wsTestData.Range("A4").Select()
